$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 34-58: renumber the "TC No" column sequentially (TC033 .. TC057),
# replacing the old placeholder "TC" values and shifting every following
# TC number up by one.
$tcNumbers = @(
    "TC033","TC034","TC035","TC036","TC037","TC038","TC039","TC040",
    "TC041","TC042","TC043","TC044","TC045","TC046","TC047","TC048",
    "TC049","TC050","TC051","TC052","TC053","TC054","TC055","TC056","TC057"
)

$row = 34
foreach ($tc in $tcNumbers) {
    $ws.Range("E$row").Value = $tc
    $row = $row + 1
}

# Row 33: fix the test-script path (the "Reports\" segment was removed).
$ws.Range("A33").Value = "Inventory\TC001CreateGoodReceipt.py"

# Widen column E slightly (closest attainable width to the target of
# 12.43 characters given this runtime's column-width granularity).
$ws.Columns("E").ColumnWidth = 11.666666666666666

# Update the view: scroll so row 13 is at the top and select K33, matching
# the author's final cursor position.
$ws.Range("K33").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
